# Adds the 20 English Premier League clubs for the 2024/2025 season to the
# "Consolidado" sheet (rows 565-584), mirroring a bulk paste of new league data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{A="Manchester City"; C="Inglaterra"; D="2024/2025"; E=20; F=1.67; G=1.1599999999999999; H=17.2; I=9.6999999999999993; J=14; K=0.86; L=66.3; M=30.15; N=10.33; O=51.13; P=61; Q=79; R=61; S=2.23},
    @{A="AFC Bournemouth"; C="Inglaterra"; D="2024/2025"; E=20; F=1.61; G=1.18; H=16.350000000000001; I=13.05; J=22; K=0.86; L=45; M=40.950000000000003; N=8.64; O=46.65; P=51; Q=73; R=46; S=2.94},
    @{A="Brighton & Hove Albion"; C="Inglaterra"; D="2024/2025"; E=20; F=1.32; G=1.27; H=14.65; I=12; J=24; K=0.88; L=45.4; M=34.5; N=8.83; O=47.6; P=57; Q=70; R=53; S=2.14},
    @{A="Nottingham Forest"; C="Inglaterra"; D="2024/2025"; E=19; F=1.1499999999999999; G=0.99; H=13.05; I=12.89; J=20; K=0.88; L=34.369999999999997; M=56.05; N=13.19; O=40.130000000000003; P=39; Q=71; R=40; S=2.2200000000000002},
    @{A="Manchester United"; C="Inglaterra"; D="2024/2025"; E=18; F=1.23; G=1.23; H=13.17; I=10.72; J=19; K=0.84; L=41.33; M=39.22; N=9.85; O=43.67; P=49; Q=65; R=54; S=1.93},
    @{A="Chelsea"; C="Inglaterra"; D="2024/2025"; E=20; F=1.72; G=1.1499999999999999; H=15.6; I=10.95; J=18; K=0.89; L=51.15; M=34.35; N=9.6199999999999992; O=46.17; P=54; Q=68; R=58; S=2.6},
    @{A="Crystal Palace"; C="Inglaterra"; D="2024/2025"; E=20; F=1.34; G=1.31; H=14.2; I=12.85; J=28; K=0.86; L=35.549999999999997; M=48.2; N=10.95; O=42.8; P=45; Q=72; R=44; S=1.77},
    @{A="Liverpool"; C="Inglaterra"; D="2024/2025"; E=18; F=1.82; G=0.84; H=16.39; I=9.61; J=17; K=0.87; L=53.83; M=32.94; N=8.9499999999999993; O=46.91; P=54; Q=71; R=57; S=2.96},
    @{A="Fulham"; C="Inglaterra"; D="2024/2025"; E=18; F=1.23; G=1.05; H=13.67; I=12.06; J=22; K=0.86; L=46.44; M=42.94; N=11.03; O=43.57; P=46; Q=65; R=52; S=2.12},
    @{A="Everton"; C="Inglaterra"; D="2024/2025"; E=19; F=0.92; G=1.27; H=11.37; I=13.32; J=23; K=0.87; L=39.049999999999997; M=46.95; N=13.09; O=42.61; P=47; Q=74; R=40; S=1.85},
    @{A="West Ham United"; C="Inglaterra"; D="2024/2025"; E=20; F=1.22; G=1.67; H=14.5; I=16.5; J=23; K=0.89; L=37.549999999999997; M=48.8; N=10.94; O=43.4; P=48; Q=71; R=46; S=1.79},
    @{A="Tottenham Hotspur"; C="Inglaterra"; D="2024/2025"; E=20; F=1.7; G=1.48; H=15; I=12.4; J=17; K=0.85; L=53; M=35.950000000000003; N=7.62; O=46.84; P=54; Q=68; R=58; S=2.11},
    @{A="Wolverhampton Wanderers"; C="Inglaterra"; D="2024/2025"; E=19; F=0.95; G=1.28; H=10.79; I=13.89; J=21; K=0.85; L=36.74; M=48.42; N=10.14; O=41.55; P=40; Q=58; R=47; S=1.43},
    @{A="Newcastle United"; C="Inglaterra"; D="2024/2025"; E=20; F=1.51; G=1.1299999999999999; H=14.25; I=13.15; J=18; K=0.88; L=40.65; M=43.45; N=9.2899999999999991; O=47.35; P=49; Q=70; R=51; S=2.4300000000000002},
    @{A="Arsenal"; C="Inglaterra"; D="2024/2025"; E=20; F=1.46; G=0.82; H=14; I=10.6; J=14; K=0.85; L=46.9; M=36.5; N=8.5399999999999991; O=48.95; P=57; Q=82; R=55; S=2.36},
    @{A="Aston Villa"; C="Inglaterra"; D="2024/2025"; E=20; F=1.39; G=1.1599999999999999; H=12.85; I=11.75; J=13; K=0.9; L=40.5; M=41.85; N=11.05; O=45.78; P=43; Q=67; R=51; S=1.51},
    @{A="Ipswich Town"; C="Inglaterra"; D="2024/2025"; E=19; F=0.84; G=1.82; H=9.68; I=15.42; J=24; K=0.83; L=32.369999999999997; M=48.84; N=13.08; O=43.85; P=49; Q=66; R=41; S=2.0099999999999998},
    @{A="Brentford"; C="Inglaterra"; D="2024/2025"; E=20; F=1.37; G=1.58; H=10.6; I=17.899999999999999; J=19; K=0.85; L=37.200000000000003; M=50.5; N=10.58; O=45.13; P=48; Q=68; R=48; S=2.15},
    @{A="Leicester City"; C="Inglaterra"; D="2024/2025"; E=20; F=0.87; G=1.66; H=8.65; I=17.600000000000001; J=26; K=0.89; L=35.15; M=50.9; N=12; O=42.16; P=42; Q=60; R=46; S=1.39},
    @{A="Southampton"; C="Inglaterra"; D="2024/2025"; E=20; F=0.86; G=2.06; H=9.3000000000000007; I=18.350000000000001; J=18; K=0.89; L=39.450000000000003; M=46.7; N=10; O=45.34; P=45; Q=60; R=52; S=1.33}
)

$startRow = 565

# Fill column A (team) for every new row first, then C (league), then D (season),
# then the numeric stat columns E..S row by row -- this reproduces the exact shared-
# string insertion order of the source edit ("Ipswich Town" added at index 293 before
# "2024/2025" at index 294).
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("A$r").Value() = $newRows[$i].A
}
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("C$r").Value() = $newRows[$i].C
}
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Range("D$r").Value() = $newRows[$i].D
}
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Range("E$r").Value() = $row.E
    $ws.Range("F$r").Value() = $row.F
    $ws.Range("G$r").Value() = $row.G
    $ws.Range("H$r").Value() = $row.H
    $ws.Range("I$r").Value() = $row.I
    $ws.Range("J$r").Value() = $row.J
    $ws.Range("K$r").Value() = $row.K
    $ws.Range("L$r").Value() = $row.L
    $ws.Range("M$r").Value() = $row.M
    $ws.Range("N$r").Value() = $row.N
    $ws.Range("O$r").Value() = $row.O
    $ws.Range("P$r").Value() = $row.P
    $ws.Range("Q$r").Value() = $row.Q
    $ws.Range("R$r").Value() = $row.R
    $ws.Range("S$r").Value() = $row.S
}

# Leave the selection where the source workbook left it after the paste.
[void]$ws.Range("F579").Select()

Write-Host "Added 20 rows (565-584) for Premier League 2024/2025"
